$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3741.6807
$ws.Range("I15").Value = 3741.6807
$ws.Range("K15").Value = 11225.0421
$ws.Range("M15").Value = -11056.0421

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 632.8570999999999
$ws.Range("I41").Value = 736
$ws.Range("J41").Value = 375
$ws.Range("K41").Value = 736
$ws.Range("L41").Value = 375
$ws.Range("M41").Value = -296
$ws.Range("N41").Value = -1255

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 6812.273
$ws.Range("I88").Value = 5846.5
$ws.Range("K88").Value = 5846.5
$ws.Range("M88").Value = -5440.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 6812.273
$ws.Range("I91").Value = 5846.5
$ws.Range("K91").Value = 5846.5
$ws.Range("M91").Value = -4442.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 543.8333
$ws.Range("I107").Value = 561.4375
$ws.Range("J107").Value = 403
$ws.Range("K107").Value = 561.4375
$ws.Range("L107").Value = 403
$ws.Range("M107").Value = 1358.5625
$ws.Range("N107").Value = -4243

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 953.4
$ws.Range("I129").Value = 197
$ws.Range("J129").Value = 1142.5
$ws.Range("K129").Value = 591
$ws.Range("L129").Value = 3427.5
$ws.Range("M129").Value = 4409
$ws.Range("N129").Value = -13427.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2676.276
$ws.Range("I132").Value = 2744.24
$ws.Range("K132").Value = 8232.719999999999
$ws.Range("M132").Value = -5702.719999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2355216.8
$ws.Range("I138").Value = 984.2820400000001
$ws.Range("J138").Value = 4351196.5
$ws.Range("K138").Value = 2952.84612
$ws.Range("L138").Value = 13053589.5
$ws.Range("M138").Value = 2187.15388
$ws.Range("N138").Value = -13063869.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1728.8948
$ws.Range("I141").Value = 1115.5625
$ws.Range("K141").Value = 3346.6875
$ws.Range("M141").Value = 1833.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5085.1978
$ws.Range("I32").Value = 3270.4265
$ws.Range("J32").Value = 14577.846
$ws.Range("K32").Value = 3270.4265
$ws.Range("L32").Value = 14577.846
$ws.Range("M32").Value = -2983.4265
$ws.Range("N32").Value = -15151.846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1491.5435
$ws.Range("I61").Value = 1360.9395
$ws.Range("J61").Value = 1823.0769
$ws.Range("K61").Value = 1360.9395
$ws.Range("L61").Value = 1823.0769
$ws.Range("M61").Value = -1148.9395
$ws.Range("N61").Value = -2247.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 372.22223
$ws.Range("I97").Value = 335.29413
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 335.29413
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 160.70587
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 30057.143
$ws.Range("J109").Value = 30057.143
$ws.Range("L109").Value = 30057.143
$ws.Range("N109").Value = -32831.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1491.5435
$ws.Range("I136").Value = 1360.9395
$ws.Range("J136").Value = 1823.0769
$ws.Range("K136").Value = 4082.8185
$ws.Range("L136").Value = 5469.2307
$ws.Range("M136").Value = -1532.8185
$ws.Range("N136").Value = -10569.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29633
$ws.Range("J51").Value = 29633
$ws.Range("L51").Value = 29633
$ws.Range("N51").Value = -30615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29079.666
$ws.Range("J55").Value = 29079.666
$ws.Range("L55").Value = 29079.666
$ws.Range("N55").Value = -29625.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1744.2307
$ws.Range("I86").Value = 1563.15
$ws.Range("K86").Value = 1563.15
$ws.Range("M86").Value = -440.1500000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1744.2307
$ws.Range("I89").Value = 1563.15
$ws.Range("K89").Value = 7815.75
$ws.Range("M89").Value = -2199.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 28000
$ws.Range("J108").Value = 28000
$ws.Range("L108").Value = 28000
$ws.Range("N108").Value = -35680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3270
$ws.Range("I113").Value = 3270
$ws.Range("K113").Value = 3270
$ws.Range("M113").Value = -1100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 648489.4399999999
$ws.Range("I134").Value = 1028751.5
$ws.Range("K134").Value = 3086254.5
$ws.Range("M134").Value = -3083719.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2164968.8
$ws.Range("I113").Value = 4329527.5
$ws.Range("J113").Value = 409.85715
$ws.Range("K113").Value = 12988582.5
$ws.Range("L113").Value = 1229.57145
$ws.Range("M113").Value = -12986412.5
$ws.Range("N113").Value = -5569.571449999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4968.5
$ws.Range("I133").Value = 2256.3333
$ws.Range("J133").Value = 6846.154
$ws.Range("K133").Value = 6768.999899999999
$ws.Range("L133").Value = 20538.462
$ws.Range("M133").Value = -1708.999899999999
$ws.Range("N133").Value = -30658.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29000
$ws.Range("J51").Value = 29000
$ws.Range("L51").Value = 29000
$ws.Range("N51").Value = -30018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23133.166
$ws.Range("J63").Value = 23460
$ws.Range("L63").Value = 23460
$ws.Range("N63").Value = -24832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 23133.166
$ws.Range("J66").Value = 23460
$ws.Range("L66").Value = 70380
$ws.Range("N66").Value = -77244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4523.109
$ws.Range("I70").Value = 4355.974
$ws.Range("J70").Value = 5454.2856
$ws.Range("K70").Value = 4355.974
$ws.Range("L70").Value = 5454.2856
$ws.Range("M70").Value = -4085.974
$ws.Range("N70").Value = -5994.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4523.109
$ws.Range("I73").Value = 4355.974
$ws.Range("J73").Value = 5454.2856
$ws.Range("K73").Value = 4355.974
$ws.Range("L73").Value = 5454.2856
$ws.Range("M73").Value = -3419.974
$ws.Range("N73").Value = -7326.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1486.5
$ws.Range("I113").Value = 1328.2727
$ws.Range("K113").Value = 1328.2727
$ws.Range("M113").Value = 841.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2515.7273
$ws.Range("I126").Value = 2980.6
$ws.Range("K126").Value = 8941.799999999999
$ws.Range("M126").Value = -6471.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2634310
$ws.Range("I132").Value = 2723.32
$ws.Range("J132").Value = 7695054
$ws.Range("K132").Value = 8169.960000000001
$ws.Range("L132").Value = 23085162
$ws.Range("M132").Value = -5639.960000000001
$ws.Range("N132").Value = -23090222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 558.0476
$ws.Range("I55").Value = 195.55556
$ws.Range("J55").Value = 829.9167
$ws.Range("K55").Value = 195.55556
$ws.Range("L55").Value = 829.9167
$ws.Range("M55").Value = -22.55556000000001
$ws.Range("N55").Value = -1175.9167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5084.769
$ws.Range("I132").Value = 5785.263
$ws.Range("J132").Value = 4419.3
$ws.Range("K132").Value = 17355.789
$ws.Range("L132").Value = 13257.9
$ws.Range("M132").Value = -14825.789
$ws.Range("N132").Value = -18317.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1889.9706
$ws.Range("I136").Value = 1418.3684
$ws.Range("K136").Value = 4255.1052
$ws.Range("M136").Value = -1705.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 56985
$ws.Range("J41").Value = 9990
$ws.Range("L41").Value = 9990
$ws.Range("N41").Value = -10770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 950
$ws.Range("I96").Value = 700
$ws.Range("J96").Value = 1075
$ws.Range("K96").Value = 700
$ws.Range("L96").Value = 1075
$ws.Range("M96").Value = 673
$ws.Range("N96").Value = -3821

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8181.5713
$ws.Range("I107").Value = 15483
$ws.Range("J107").Value = 880.1429000000001
$ws.Range("K107").Value = 46449
$ws.Range("L107").Value = 2640.4287
$ws.Range("M107").Value = -44529
$ws.Range("N107").Value = -6480.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 61075.8
$ws.Range("I122").Value = 13634.267
$ws.Range("J122").Value = 203400.4
$ws.Range("K122").Value = 40902.801
$ws.Range("L122").Value = 610201.2
$ws.Range("M122").Value = -38452.801
$ws.Range("N122").Value = -615101.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 100000820
$ws.Range("I126").Value = 812
$ws.Range("J126").Value = 333334180
$ws.Range("K126").Value = 2436
$ws.Range("L126").Value = 1000002540
$ws.Range("M126").Value = 34
$ws.Range("N126").Value = -1000007480

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3194.75
$ws.Range("I132").Value = 3814.8635
$ws.Range("K132").Value = 11444.5905
$ws.Range("M132").Value = -8914.5905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2109.3276
$ws.Range("I136").Value = 2032.7435
$ws.Range("J136").Value = 2266.5264
$ws.Range("K136").Value = 6098.2305
$ws.Range("L136").Value = 6799.5792
$ws.Range("M136").Value = -3548.2305
$ws.Range("N136").Value = -11899.5792
